# Capstone_DataScience_Dhinakaran.pptx edit
# 1) Add a "Challenges" entry to the Content (TOC) slide, right before "Results".
# 2) Insert a brand-new "Challenges" slide (position 6, right after "Workflow"
#    and before "Results"), duplicating the "Conclusion" slide's layout/shape
#    structure (Title + text Content Placeholder + Slide Number placeholder)
#    and filling in the challenges text.

$p = $ppt.ActivePresentation

# --- 1. TOC slide: insert "Challenges" paragraph before "Results " ---
$toc = $p.Slides.Item(2)
$tocBody = $toc.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $tocBody.Paragraphs().Count; $i++) {
    if ($tocBody.Paragraphs($i).Text -eq "Results ") {
        $tocBody.Paragraphs($i).InsertBefore("Challenges`r")
        break
    }
}

# --- 2. New "Challenges" slide ---
# Duplicate the "Conclusion" slide (slide 7) since it already has the same
# shape layout we need (Title + text Content Placeholder + Slide Number).
$conclusion = $p.Slides.Item(7)
$dupRange = $conclusion.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(6)

$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Challenges"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "FourSquare result will have all information related to Indian and other restaurant, like Indian consulate, Indian Yoga center, Chinese restaurants etc, because API is designed to split the search string " + [char]8220 + "Indian Restaurant" + [char]8221 + " like " + [char]8220 + "Indian" + [char]8221 + ", " + [char]8220 + "Restaurant" + [char]8221 + ", " + [char]8220 + "Indian Restaurant" + [char]8221 + " and " + [char]8220 + "Indian Restaurants" + [char]8221 + " and perform the search. So, analyzing the result and extract correct the information is vital to this project.`r" + `
"And search area is important to this analysis, to avoid overlapping of information which will mislead the results. So, utilized different data (ward area information from UK gov site) to calculate average ward area. Total 625 wards are in London (including inner and outer London region) with total 1594 sq. meter " + [char]8211 + " which make search limit average to 2500 sq. meter.`r" + `
"Because of some generic location name used in both USA and UK, google API returned US coordinates, so concatenated the Ward, Borough and country to make unique location and fetched correct and accurate coordinates from google API.`r"
